$d = $word.ActiveDocument

# Whole-paragraph (single run) text replacements — assigning directly to
# Paragraph.Range.Text lets Word recompute xml:space the same way the
# original authoring tool did.
function Set-ParaText($matchText, $newText) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]7, [char]13)
        if ($t -eq $matchText) {
            $p.Range.Text = $newText
            return $true
        }
    }
    return $false
}

$replacements = @(
    @("Dna. Ioana Mirea", "Dl. 111 111"),
    @("Şef Birou Aprovizionare", "111"),
    @("Birou Achiziţii", "111"),
    @("S.C. Automatica S.A.", "S.C. Apex Electric VD S.R.L."),
    @("Șoseaua Afumați Nr. 108, Voluntari", "Strada Izvorul Rece Nr. 3, Bl. 14, Ap. 3, Brașov"),
    @("Tel: +40 (372) 058 100*107", "Tel: 111"),
    @("Fax: +40 (372) 058 101", "Fax: 111"),
    @("Mobil: +40 (729) 035 164", "Mobil: 111")
)

foreach ($pair in $replacements) {
    $ok = Set-ParaText $pair[0] $pair[1]
    if (-not $ok) {
        Write-Output "WARNING: paragraph text not found for replacement: $($pair[0])"
    }
}

# The E-mail line has two runs ("E-mail: " and the hyperlink-styled
# address) — only the address run changes, so locate it with Find and
# replace just that span via an independently constructed Range (avoids
# reusing the paragraph Range object for the write).
$r = $d.Content
$found = $r.Find.Execute("ioana.mirea@automatica.ro", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sub = $d.Range($r.Start, $r.End)
    $sub.Text = "111@111.com"
} else {
    Write-Output "WARNING: email address run not found"
}
